$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25
$ws.Range("A25").Value = "BonusPower"
$ws.Range("B25").Value = 0.02
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 24
$ws.Range("E25").Value = 500
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = "lose"

# Row 26
$ws.Range("A26").Value = "BonusPower"
$ws.Range("B26").Value = 2
$ws.Range("C26").Value = 1070
$ws.Range("D26").Value = 127
$ws.Range("E26").Value = 1000
$ws.Range("F26").Value = 2000
$ws.Range("G26").Value = "win"
